$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2,3,4) have their per-row varying values (Fecha/D, Volumen/J,
# Precio minimo/K, Precio maximo/L, Precio promedio ponderado/M, Precio $/Kg/P) cyclically
# shifted up by one row: row2 <- row3, row3 <- row4, row4 <- row2 (wrap-around).
# Capture originals first using Value2 (Value is unreliable in this runtime).
$D2 = $ws.Range("D2").Value2
$J2 = $ws.Range("J2").Value2
$K2 = $ws.Range("K2").Value2
$L2 = $ws.Range("L2").Value2
$M2 = $ws.Range("M2").Value2
$P2 = $ws.Range("P2").Value2

$D3 = $ws.Range("D3").Value2
$J3 = $ws.Range("J3").Value2
$K3 = $ws.Range("K3").Value2
$L3 = $ws.Range("L3").Value2
$M3 = $ws.Range("M3").Value2
$P3 = $ws.Range("P3").Value2

$D4 = $ws.Range("D4").Value2
$J4 = $ws.Range("J4").Value2
$K4 = $ws.Range("K4").Value2
$L4 = $ws.Range("L4").Value2
$M4 = $ws.Range("M4").Value2
$P4 = $ws.Range("P4").Value2

# Row 2 gets row 3's original values
$ws.Range("D2").Value2 = $D3
$ws.Range("J2").Value2 = $J3
$ws.Range("K2").Value2 = $K3
$ws.Range("L2").Value2 = $L3
$ws.Range("M2").Value2 = $M3
$ws.Range("P2").Value2 = $P3

# Row 3 gets row 4's original values
$ws.Range("D3").Value2 = $D4
$ws.Range("J3").Value2 = $J4
$ws.Range("K3").Value2 = $K4
$ws.Range("L3").Value2 = $L4
$ws.Range("M3").Value2 = $M4
$ws.Range("P3").Value2 = $P4

# Row 4 gets row 2's original values
$ws.Range("D4").Value2 = $D2
$ws.Range("J4").Value2 = $J2
$ws.Range("K4").Value2 = $K2
$ws.Range("L4").Value2 = $L2
$ws.Range("M4").Value2 = $M2
$ws.Range("P4").Value2 = $P2
